$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3580396666666667
$ws.Cells.Item(2, 8).Value = 1.074119
$ws.Cells.Item(2, 9).Value = 0.07478320321406828
$ws.Cells.Item(2, 10).Value = 0.07478320321406827
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3580396666666667
$ws.Cells.Item(2, 14).Value = 1.074119
$ws.Cells.Item(2, 15).Value = 0.07478320321406828
$ws.Cells.Item(2, 16).Value = 0.07478320321406827
$ws.Cells.Item(2, 17).Value = 0.1281924029067778
$ws.Cells.Item(2, 18).Value = 1.153731626161
$ws.Cells.Item(2, 19).Value = 0.005592527482956633
$ws.Cells.Item(2, 20).Value = 0.00559252748295663

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3580396666666667
$ws.Cells.Item(3, 8).Value = 1.074119
$ws.Cells.Item(3, 9).Value = 0.07478320321406828
$ws.Cells.Item(3, 10).Value = 0.07478320321406827
$ws.Cells.Item(3, 15).Value = 0.555484973478924
$ws.Cells.Item(3, 16).Value = 0.555484973478924
$ws.Cells.Item(3, 17).Value = 0.952205180687889
$ws.Cells.Item(3, 18).Value = 8.569846626191
$ws.Cells.Item(3, 19).Value = 0.0415409456540357
$ws.Cells.Item(3, 20).Value = 0.0415409456540357

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3580396666666667
$ws.Cells.Item(4, 8).Value = 1.074119
$ws.Cells.Item(4, 9).Value = 0.07478320321406828
$ws.Cells.Item(4, 10).Value = 0.07478320321406827
$ws.Cells.Item(4, 15).Value = 0.3697318233070078
$ws.Cells.Item(4, 16).Value = 0.3697318233070078
$ws.Cells.Item(4, 17).Value = 0.6337895252381112
$ws.Cells.Item(4, 18).Value = 5.704105727143
$ws.Cells.Item(4, 19).Value = 0.02764973007707595
$ws.Cells.Item(4, 20).Value = 0.02764973007707594

# Row 5
$ws.Cells.Item(5, 9).Value = 0.555484973478924
$ws.Cells.Item(5, 10).Value = 0.555484973478924
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3580396666666667
$ws.Cells.Item(5, 14).Value = 1.074119
$ws.Cells.Item(5, 15).Value = 0.07478320321406828
$ws.Cells.Item(5, 16).Value = 0.07478320321406827
$ws.Cells.Item(5, 17).Value = 0.952205180687889
$ws.Cells.Item(5, 18).Value = 8.569846626191
$ws.Cells.Item(5, 19).Value = 0.0415409456540357
$ws.Cells.Item(5, 20).Value = 0.0415409456540357

# Row 6
$ws.Cells.Item(6, 9).Value = 0.555484973478924
$ws.Cells.Item(6, 10).Value = 0.555484973478924
$ws.Cells.Item(6, 15).Value = 0.555484973478924
$ws.Cells.Item(6, 16).Value = 0.555484973478924
$ws.Cells.Item(6, 19).Value = 0.3085635557608809
$ws.Cells.Item(6, 20).Value = 0.3085635557608809

# Row 7
$ws.Cells.Item(7, 9).Value = 0.555484973478924
$ws.Cells.Item(7, 10).Value = 0.555484973478924
$ws.Cells.Item(7, 15).Value = 0.3697318233070078
$ws.Cells.Item(7, 16).Value = 0.3697318233070078
$ws.Cells.Item(7, 19).Value = 0.2053804720640074
$ws.Cells.Item(7, 20).Value = 0.2053804720640074

# Row 8
$ws.Cells.Item(8, 9).Value = 0.3697318233070078
$ws.Cells.Item(8, 10).Value = 0.3697318233070078
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.3580396666666667
$ws.Cells.Item(8, 14).Value = 1.074119
$ws.Cells.Item(8, 15).Value = 0.07478320321406828
$ws.Cells.Item(8, 16).Value = 0.07478320321406827
$ws.Cells.Item(8, 17).Value = 0.6337895252381112
$ws.Cells.Item(8, 18).Value = 5.704105727143
$ws.Cells.Item(8, 19).Value = 0.02764973007707595
$ws.Cells.Item(8, 20).Value = 0.02764973007707594

# Row 9
$ws.Cells.Item(9, 9).Value = 0.3697318233070078
$ws.Cells.Item(9, 10).Value = 0.3697318233070078
$ws.Cells.Item(9, 15).Value = 0.555484973478924
$ws.Cells.Item(9, 16).Value = 0.555484973478924
$ws.Cells.Item(9, 19).Value = 0.2053804720640074
$ws.Cells.Item(9, 20).Value = 0.2053804720640074

# Row 10
$ws.Cells.Item(10, 9).Value = 0.3697318233070078
$ws.Cells.Item(10, 10).Value = 0.3697318233070078
$ws.Cells.Item(10, 15).Value = 0.3697318233070078
$ws.Cells.Item(10, 16).Value = 0.3697318233070078
$ws.Cells.Item(10, 19).Value = 0.1367016211659245
$ws.Cells.Item(10, 20).Value = 0.1367016211659244
